# Updates the grades sheet with recomputed scores (PRD detection fix /
# non-standard filename handling fix). Raw Score (B), Max Possible (C),
# Percentage (D) and Grade (F) move for nearly every student because the
# evaluation max shifted from 20.644 -> 20.754; rows 14-15 additionally
# gained a Criteria Count point (H: 14 -> 15) from the fixed PRD match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> @{ col = value }
$updates = @{
    2  = @{ B = 19.794; C = 20.754; D = 95.37438566059554; F = 96.37438566059554 }
    3  = @{ B = 19.794; C = 20.754; D = 95.37438566059554; F = 96.37438566059554 }
    4  = @{ B = 19.794; C = 20.754; D = 95.37438566059554; F = 96.37438566059554 }
    5  = @{ B = 18.329; C = 20.754; D = 88.31550544473355; F = 89.31550544473355 }
    6  = @{ B = 18.329; C = 20.754; D = 88.31550544473355; F = 89.31550544473355 }
    7  = @{ B = 18.069; C = 20.754; D = 87.06273489447817; F = 87.06273489447817 }
    8  = @{ B = 18.069; C = 20.754; D = 87.06273489447817; F = 87.06273489447817 }
    9  = @{ B = 17.139; C = 20.754; D = 82.58167100318011; F = 83.58167100318011 }
    10 = @{ B = 17.139; C = 20.754; D = 82.58167100318011; F = 83.58167100318011 }
    11 = @{ B = 15.662; C = 20.754; D = 75.46497060807555; F = 75.46497060807555 }
    12 = @{ B = 15.662; C = 20.754; D = 75.46497060807555; F = 75.46497060807555 }
    13 = @{            C = 20.754; D = 74.27002023706272; F = 74.27002023706272 }
    14 = @{ B = 15.167; C = 20.754; D = 73.07988821432012; F = 73.07988821432012; H = 15 }
    15 = @{ B = 15.167; C = 20.754; D = 73.07988821432012; F = 73.07988821432012; H = 15 }
    16 = @{ B = 13.712; C = 20.754; D = 66.06919148116025; F = 66.06919148116025 }
    17 = @{ B = 13.712; C = 20.754; D = 66.06919148116025; F = 66.06919148116025 }
    18 = @{ B = 12.719; C = 20.754; D = 61.28457164883877; F = 61.28457164883877 }
    19 = @{ B = 12.182; C = 20.754; D = 58.69711862773441; F = 58.69711862773441 }
    20 = @{ B = 12.182; C = 20.754; D = 58.69711862773441; F = 58.69711862773441 }
    21 = @{ B = 7.135;  C = 20.754; D = 34.37891490796954; F = 34.37891490796954 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
